# Updates the cost parameter values in the Investment_Cost_Overview workbook.
$wb = $excel.ActiveWorkbook

$wsCost = $wb.Worksheets.Item("Investment_Cost")
$wsSources = $wb.Worksheets.Item("Sources")

# --- Row 7 (Methanol_Plant): cost is a placeholder "1" in every year column ---
$wsCost.Range("B7:F7").Value = 1

# --- Row 9 (Methanol_storage): refreshed cost figures ---
$wsCost.Range("B9:F9").Value = 139.58682300390799

# --- Row 10 (Hydrogen_storage): values re-based from fractional (0.xxx) to
#     absolute Euro figures, switching the number format accordingly ---
$wsCost.Range("B10:F10").NumberFormat = "#,##0.00"
$wsCost.Range("B10").Value = 121000
$wsCost.Range("C10").Value = 170500
$wsCost.Range("D10").Value = 99000
$wsCost.Range("E10").Value = 61000
$wsCost.Range("F10").Value = 46000

# --- Stray formatted-but-empty cells left behind below the table ---
$wsCost.Range("B16").NumberFormat = "#,##0.00"
$wsCost.Range("E26").NumberFormat = "#,##0"

# --- Note explaining the Methanol_Plant placeholder cost, added to Sources ---
$wsSources.Range("A3").Value = "The methanol plant value is set to 1 to give the model a cost. The investment itself is in the distillation tower for the output of methanol. "

# --- Restore on-screen selections ---
$wsCost.Activate() | Out-Null
$wsCost.Range("I24").Select() | Out-Null
$wsSources.Activate() | Out-Null
$wsSources.Range("D6").Select() | Out-Null
$wsCost.Activate() | Out-Null
